$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.969.88'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.67%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.311.81'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.92%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.81'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.73'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.45%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.99%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.310.88'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.95%  '

$ws.Range('E10').Value = '  -0.65%  '

$ws.Range('E11').Value = '  -0.52%  '

$ws.Range('E12').Value = '  +0.77%  '

$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.92'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.16%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.725.19'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.90%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.848.27'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.49%  '

$ws.Range('E17').Value = '  +0.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.320.13'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.61%  '

$ws.Range('E20').Value = '  -2.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.57'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('E22').Value = '  +2.77%  '

$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.38'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.174'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.08%  '

$ws.Range('E26').Value = '  +0.17%  '

$ws.Range('E27').Value = '  -1.51%  '

$ws.Range('E28').Value = '  -0.79%  '

$ws.Range('E29').Value = '  -0.61%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.10'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.40%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0740'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.34%  '

$ws.Range('E32').Value = '  +4.88%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.90'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.35%  '

$ws.Range('E34').Value = '  +2.82%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.34'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +7.29%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.00%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.91'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.25%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.14'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.74%  '

$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '307.43'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.23%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.52'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '141.60'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.22%  '

$ws.Range('E43').Value = '  +1.14%  '

$ws.Range('E44').Value = '  +1.27%  '

$ws.Range('E45').Value = '  -0.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.558'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.40'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0212'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.08%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.01'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.53'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.80%  '
